$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Refresh the "panel_query_time" style timestamps on the "data"
#    sheet (column F, rows 2-31) to the newer query time.
# ------------------------------------------------------------------
$dataSheet = $wb.Worksheets.Item("data")

$dataSheet.Range("F2").Value = "2021-10-05 14:21:07.607854"
$dataSheet.Range("F3").Value = "2021-10-05 14:21:07.607861"
$dataSheet.Range("F4").Value = "2021-10-05 14:21:07.607865"
$dataSheet.Range("F5").Value = "2021-10-05 14:21:07.607867"
$dataSheet.Range("F6").Value = "2021-10-05 14:21:07.607870"
$dataSheet.Range("F7").Value = "2021-10-05 14:21:07.607873"
$dataSheet.Range("F8").Value = "2021-10-05 14:21:07.607876"
$dataSheet.Range("F9").Value = "2021-10-05 14:21:07.607878"
$dataSheet.Range("F10").Value = "2021-10-05 14:21:07.607881"
$dataSheet.Range("F11").Value = "2021-10-05 14:21:07.607884"
$dataSheet.Range("F12").Value = "2021-10-05 14:21:07.607887"
$dataSheet.Range("F13").Value = "2021-10-05 14:21:07.607889"
$dataSheet.Range("F14").Value = "2021-10-05 14:21:07.607892"
$dataSheet.Range("F15").Value = "2021-10-05 14:21:07.607895"
$dataSheet.Range("F16").Value = "2021-10-05 14:21:07.607897"
$dataSheet.Range("F17").Value = "2021-10-05 14:21:07.607900"
$dataSheet.Range("F18").Value = "2021-10-05 14:21:07.607902"
$dataSheet.Range("F19").Value = "2021-10-05 14:21:07.607905"
$dataSheet.Range("F20").Value = "2021-10-05 14:21:07.607908"
$dataSheet.Range("F21").Value = "2021-10-05 14:21:07.607910"
$dataSheet.Range("F22").Value = "2021-10-05 14:21:07.607913"
$dataSheet.Range("F23").Value = "2021-10-05 14:21:07.607915"
$dataSheet.Range("F24").Value = "2021-10-05 14:21:07.607918"
$dataSheet.Range("F25").Value = "2021-10-05 14:21:07.607920"
$dataSheet.Range("F26").Value = "2021-10-05 14:21:07.607923"
$dataSheet.Range("F27").Value = "2021-10-05 14:21:07.607926"
$dataSheet.Range("F28").Value = "2021-10-05 14:21:07.607928"
$dataSheet.Range("F29").Value = "2021-10-05 14:21:07.607931"
$dataSheet.Range("F30").Value = "2021-10-05 14:21:07.607933"
$dataSheet.Range("F31").Value = "2021-10-05 14:21:07.607936"

# ------------------------------------------------------------------
# 2. Add a new "metadata" tab (placed right after "data") describing
#    the panel query that produced this workbook.
# ------------------------------------------------------------------
$newWs = $wb.Worksheets.Add()
$newWs.Name = "metadata"

# Re-resolve sheets by name - after adding/renaming a sheet, older
# variables can end up tracking stale positions rather than the sheet
# itself.
$dataSheet = $wb.Worksheets.Item("data")
$metaSheet = $wb.Worksheets.Item("metadata")

# Pull the header style (bold/bordered/centered) from the "data" sheet
# header row so "metadata" matches the same look - copying also tiles
# the 5-wide source across the 6-wide destination range.
$headerSrc = $dataSheet.Range("B1:F1")
$headerDst = $metaSheet.Range("B1:G1")
$headerSrc.Copy($headerDst)

# Same for the A2 "index" style.
$idxSrc = $dataSheet.Range("A2")
$idxDst = $metaSheet.Range("A2")
$idxSrc.Copy($idxDst)

# Header labels
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Inherited non-medullary thyroid cancer"
$metaSheet.Range("C2").Value = 171

# "data_version" must be stored as text ("1.5"), not a number.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.5"
$metaSheet.Range("D2").Style = "Normal"

$metaSheet.Range("E2").Value = "2020-12-02T12:42:07.841648Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:21:07.604600"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/171/?format=json"

# ------------------------------------------------------------------
# 3. Put "metadata" right after "data" in the tab order.
# ------------------------------------------------------------------
$dataSheet = $wb.Worksheets.Item("data")
$metaSheet = $wb.Worksheets.Item("metadata")
$metaSheet.Move($null, $dataSheet)

Write-Output "metadata tab added"
